$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "DeathPercentage" column (C): shift "total_vaccinations" data
# from column D into column C, then clear the now-unused column D.
$ws.Range("C1").Value = "total_vaccinations"
$ws.Range("C2").Value = 10701388790
$ws.Range("D1:D2").Clear()

# Update the sheet's selection to reflect the new used range (A1:C2).
$ws.Range("A1:C2").Select()
